# Update symbol list (crypto price/coin data refresh) on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price column (D) updates for existing rows ---------------------------
# These values look numeric, so force text format first to keep the exact
# original string representation (trailing/leading zeros, precision) instead
# of letting Excel coerce them into floating point numbers.
$priceUpdates = @{
    "D2"  = "244.17"
    "D3"  = "23.98"
    "D4"  = "5.257"
    "D5"  = "0.05834"
    "D6"  = "6.459"
    "D7"  = "3.236"
    "D8"  = "0.8083"
    "D9"  = "0.8884"
    "D10" = "0.1382"
    "D11" = "0.07110"
    "D12" = "0.03086"
    "D13" = "0.03039"
    "D14" = "0.09333"
    "D15" = "3.827"
    "D16" = "0.001534"
    "D17" = "0.04713"
    "D18" = "0.0006040"
    "D19" = "0.006178"
    "D21" = "0.004071"
    "D22" = "0.00008701"
    "D24" = "2.168"
    "D25" = "0.3185"
    "D28" = "0.0002329"
    "D40" = "0.03841"
    "D41" = "0.1054"
    "D42" = "0.002557"
    "D43" = "0.006291"
    "D44" = "0.007291"
    "D45" = "0.00005318"
    "D47" = "0.5213"
    "D48" = "0.002143"
}

foreach ($addr in $priceUpdates.Keys) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $priceUpdates[$addr]
}

# --- Rows 41-43 got reordered / replaced with new coin entries ------------
# Row 41: KickToken -> BKEXToken
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("E41").Value = "40BKEXTokenBKK"

# Row 42: BKEXToken -> CEJI
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("E42").Value = "41CEJICEJIWorstin24h"

# Row 43: CEJI -> KickToken
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("E43").Value = "42KickTokenKICK"

# Row 47: label lost its "Worstin24h" suffix
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
